$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the contents of two same-shaped cells (column B through AC)
# between two rows. Column A (the running index) is intentionally left
# untouched by the caller - these row pairs only swap their data payload.
# ---------------------------------------------------------------------------
function Swap-Cell($sheet, $addr1, $addr2) {
    $v1 = $sheet.Range($addr1).Value2
    $v2 = $sheet.Range($addr2).Value2
    $sheet.Range($addr1).Value = $v2
    $sheet.Range($addr2).Value = $v1
}

function Swap-Rows($sheet, $row1, $row2) {
    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    foreach ($col in $cols) {
        $a1 = $col + $row1
        $a2 = $col + $row2
        Swap-Cell $sheet $a1 $a2
    }
}

# Updated-data swaps (re-ordering within the source feed moved these match
# records across the two adjacent rows they already occupied).
Swap-Rows $ws 27 28
Swap-Rows $ws 104 105
Swap-Rows $ws 107 108
Swap-Rows $ws 142 145
Swap-Rows $ws 148 150

# ---------------------------------------------------------------------------
# Row 173 gets new match data (it used to be the placeholder/incomplete last
# row), and two brand-new rows (174, 175) are appended after it.
# ---------------------------------------------------------------------------

# Give the new rows 174/175 the same formatting as row 173 (bold/bordered
# index cell in column A, date-formatted column E) before filling values in,
# so the cloned style entries line up with the existing style table instead
# of minting new ones.
$ws.Range("A173").Copy()
$ws.Range("A174:A175").PasteSpecial(-4122)
$ws.Range("E173").Copy()
$ws.Range("E174:E175").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 173 (rewritten in place)
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 7801327
$ws.Range("C173").Value = "Bolivia Primera División"
$ws.Range("D173").Value = "Bolivia Apertura"
$ws.Range("E173").Value = 45348.875
$ws.Range("F173").Value = "Real Tomayapo"
$ws.Range("G173").Value = "Royal Pari FC"
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 2
$ws.Range("J173").Value = "A"
$ws.Range("K173").Value = 1.666
$ws.Range("L173").Value = 3.6
$ws.Range("M173").Value = 4.5
$ws.Range("N173").Value = 1.85
$ws.Range("O173").Value = 3.5
$ws.Range("P173").Value = 4.2
$ws.Range("Q173").Value = -0.5
$ws.Range("R173").Value = 1.85
$ws.Range("S173").Value = 1.95
$ws.Range("T173").Value = 2.5
$ws.Range("U173").Value = 1.85
$ws.Range("V173").Value = 1.95
$ws.Range("W173").Value = -1
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = 3.2
$ws.Range("Z173").Value = -1
$ws.Range("AA173").Value = 0.95
$ws.Range("AB173").Value = -1
$ws.Range("AC173").Value = 0.95

# Row 174 (new)
$ws.Range("A174").Value = 172
$ws.Range("B174").Value = 7801328
$ws.Range("C174").Value = "Bolivia Primera División"
$ws.Range("D174").Value = "Bolivia Apertura"
$ws.Range("E174").Value = 45349.875
$ws.Range("F174").Value = "Independiente Petrolero"
$ws.Range("G174").Value = "Nacional Potosi"
$ws.Range("H174").Value = 1
$ws.Range("I174").Value = 1
$ws.Range("J174").Value = "D"
$ws.Range("K174").Value = 2.5
$ws.Range("L174").Value = 3.4
$ws.Range("M174").Value = 2.6
$ws.Range("N174").Value = 2.4
$ws.Range("O174").Value = 3.4
$ws.Range("P174").Value = 2.8
$ws.Range("Q174").Value = 0
$ws.Range("R174").Value = 1.75
$ws.Range("S174").Value = 2.05
$ws.Range("T174").Value = 2.75
$ws.Range("U174").Value = 1.9
$ws.Range("V174").Value = 1.9
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 2.4
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = 0
$ws.Range("AA174").Value = -0.0
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = 0.8999999999999999

# Row 175 (new)
$ws.Range("A175").Value = 173
$ws.Range("B175").Value = 7892752
$ws.Range("C175").Value = "Bolivia Primera División"
$ws.Range("D175").Value = "Bolivia Apertura"
$ws.Range("E175").Value = 45350.66666666666
$ws.Range("F175").Value = "San Antonio Bulo Bulo"
$ws.Range("G175").Value = "The Strongest"
$ws.Range("H175").Value = 5
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = "H"
$ws.Range("K175").Value = 3.3
$ws.Range("L175").Value = 3.5
$ws.Range("M175").Value = 1.909
$ws.Range("N175").Value = 3.3
$ws.Range("O175").Value = 3.6
$ws.Range("P175").Value = 1.85
$ws.Range("Q175").Value = 0.5
$ws.Range("R175").Value = 1.8
$ws.Range("S175").Value = 2
$ws.Range("T175").Value = 2.75
$ws.Range("U175").Value = 1.95
$ws.Range("V175").Value = 1.85
$ws.Range("W175").Value = 2.3
$ws.Range("X175").Value = -1
$ws.Range("Y175").Value = -1
$ws.Range("Z175").Value = 0.8
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = 0.95
$ws.Range("AC175").Value = -1
